$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update container height and diameter (variable cost driver tied to demand model)
$ws.Range("C4").Value = 0.58874952899999999
$ws.Range("C5").Value = 0.58874952899999999
